$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in G2: the sentence was missing the "'스페이스 바'" phrase.
$ws.Range("G2").Value = "세션 4가 끝났습니다.`n`n쉬는 시간입니다.`n`n충분히 휴식을 취한 뒤, '스페이스 바' 를 눌러 다음으로 진행해주세요."

# Update the saved view/selection of the sheet: scroll the window up by one
# row (topLeftCell A6 -> A5) and move the active selection to L5.
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
$ws.Range("L5").Select()
$excel.Goto($ws.Range("L5"), $false)
